$d = $word.ActiveDocument

# ------------------------------------------------------------------
# RF4.6 - Gestione Ordine Cliente
#   " dovrà poter accedere ad uno dei ordine che ha preso in carico."
#   -> " dovrà poter accedere ad uno degli ordini a lui associati."
# ------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $d.Paragraphs($i)
    $t = $para.Range.Text
    if ($t -like "*ad uno dei ordine che ha preso in carico.*") {
        $para.Range.Find.Execute(
            "ad uno dei ordine che ha preso in carico.",
            $true, $false, $false, $false, $false,
            $true, 1, $false,
            "ad uno degli ordini a lui associati.",
            2
        )
        break
    }
}

# ------------------------------------------------------------------
# RF4.7 - Gestione Preventivo Cliente (first occurrence only - the
# near-identical RF5.6 paragraph must stay untouched)
#   " dovrà poter accedere ad uno dei preventivi che ha preso in carico."
#   -> " dovrà poter accedere ad uno dei preventivi a lui associati."
# ------------------------------------------------------------------
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $d.Paragraphs($i)
    $t = $para.Range.Text
    if ($t -like "*ad uno dei preventivi che ha preso in carico.*") {
        $para.Range.Find.Execute(
            "ad uno dei preventivi che ha preso in carico.",
            $true, $false, $false, $false, $false,
            $true, 1, $false,
            "ad uno dei preventivi a lui associati.",
            2
        )
        break
    }
}
